$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the logged test-run values in row 2 ---
# A2: usuarioAp value used for the run
$ws.Range("A2").Value = "jtangt"

# C2: account number (must stay text, keep leading format as text like before)
$ws.Range("C2").Value = "'1001678945"

# G2: transaction id returned by the run
$ws.Range("G2").Value = "AAACT2318450ZF4MF"

# H2: timestamp of the run completion
$ws.Range("H2").Value = "3 jul. 2023, 11:22:05"

# --- Wrap text on A2 (new cell style used for the long usuarioAp value) ---
$ws.Range("A2").WrapText = $true

# --- Widen column G to fit the new transaction id ---
$ws.Columns("G").ColumnWidth = 18.529947916666668

# --- Update the active selection left by the last editor ---
[void]$ws.Range("I10").Select()
